$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column L (12th column) from stored OOXML width 13 to stored width 11.
# Excel COM ColumnWidth is offset from the stored OOXML <col width> by ~5/6
# (the default column padding), so use 11 - 5/6 to land exactly on 11 once saved.
$ws.Columns.Item(12).ColumnWidth = 11 - (5/6)

# Clear the text content of the cells that were blanked out in row 2,
# leaving them as empty cells (matching the already-empty D2/E2/N2 cells).
$ws.Range("C2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("P2").Value = ""
